# Generate Report for Handback
# Adds a third handback row (ca143e9c-70a9-40bc-a684-5517b1eb4cf1) to the
# "Overview", "zh-cn" and "de-de" sheets, growing each sheet's table by one
# row and wiring up the matching hyperlinks.

$wb = $excel.ActiveWorkbook

$fileName      = "ca143e9c-70a9-40bc-a684-5517b1eb4cf1.md"
$pathAndName   = "e2e\ca143e9c-70a9-40bc-a684-5517b1eb4cf1.md"
$extension     = ".md"
$statusInSync  = "Handed back: in sync with en-US"
$genDate       = "2016-08-19 22:51:01"

$xliffZh       = "ca143e9c-70a9-40bc-a684-5517b1eb4cf1.83bf196a4a2a17b09863015afd929f9aa74d0463.zh-cn.xlf"
$xliffZhHoDate = "2016-08-19 22:50:56"
$xliffZhHbDate = "2016-08-19 22:51:27"

$xliffDe       = "ca143e9c-70a9-40bc-a684-5517b1eb4cf1.83bf196a4a2a17b09863015afd929f9aa74d0463.de-de.xlf"
$xliffDeHoDate = "2016-08-19 22:51:01"
$xliffDeHbDate = "2016-08-19 22:51:34"

# ---------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$loOverview = $wsOverview.ListObjects.Item(1)
$loOverview.ListRows.Add() | Out-Null

$wsOverview.Range("A4").Value = $fileName
$wsOverview.Range("B4").Value = $pathAndName
$wsOverview.Range("B4").Style = "HyperLink"
$wsOverview.Range("C4").Value = $extension
$wsOverview.Range("E4").Value = $statusInSync
$wsOverview.Range("F4").Value = $statusInSync
$wsOverview.Range("G4").Value = $genDate
$wsOverview.Range("G4").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$wsOverview.Hyperlinks.Add(
    $wsOverview.Range("B4"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/83bf196a4a2a17b09863015afd929f9aa74d0463/e2e/ca143e9c-70a9-40bc-a684-5517b1eb4cf1.md",
    "",
    "",
    $pathAndName
) | Out-Null

# ---------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")
$loZh = $wsZh.ListObjects.Item(1)
$loZh.ListRows.Add() | Out-Null

$wsZh.Range("A4").Value = $fileName
$wsZh.Range("A4").Style = "HyperLink"
$wsZh.Range("B4").Value = $extension
$wsZh.Range("C4").Value = $statusInSync
$wsZh.Range("D4").Value = "e2e"
$wsZh.Range("E4").Value = "ht"
$wsZh.Range("F4").Value = "'True"
$wsZh.Range("G4").Value = $xliffZh
$wsZh.Range("H4").Value = $xliffZhHoDate
$wsZh.Range("H4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZh.Range("I4").Value = $fileName
$wsZh.Range("I4").Style = "HyperLink"
$wsZh.Range("J4").Value = $xliffZh
$wsZh.Range("K4").Value = $xliffZhHbDate
$wsZh.Range("K4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZh.Range("L4").Value = "'"
$wsZh.Range("M4").Value = "'True"
$wsZh.Range("N4").Value = "'"
$wsZh.Range("O4").Value = "'False"
$wsZh.Range("P4").Value = "'"

$wsZh.Hyperlinks.Add(
    $wsZh.Range("A4"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/83bf196a4a2a17b09863015afd929f9aa74d0463/e2e/ca143e9c-70a9-40bc-a684-5517b1eb4cf1.md",
    "",
    "",
    $fileName
) | Out-Null
$wsZh.Hyperlinks.Add(
    $wsZh.Range("I4"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/83bf196a4a2a17b09863015afd929f9aa74d0463/e2e/ca143e9c-70a9-40bc-a684-5517b1eb4cf1.md",
    "",
    "",
    $fileName
) | Out-Null

# ---------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")
$loDe = $wsDe.ListObjects.Item(1)
$loDe.ListRows.Add() | Out-Null

$wsDe.Range("A4").Value = $fileName
$wsDe.Range("A4").Style = "HyperLink"
$wsDe.Range("B4").Value = $extension
$wsDe.Range("C4").Value = $statusInSync
$wsDe.Range("D4").Value = "e2e"
$wsDe.Range("E4").Value = "ht"
$wsDe.Range("F4").Value = "'True"
$wsDe.Range("G4").Value = $xliffDe
$wsDe.Range("H4").Value = $xliffDeHoDate
$wsDe.Range("H4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDe.Range("I4").Value = $fileName
$wsDe.Range("I4").Style = "HyperLink"
$wsDe.Range("J4").Value = $xliffDe
$wsDe.Range("K4").Value = $xliffDeHbDate
$wsDe.Range("K4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDe.Range("L4").Value = "'"
$wsDe.Range("M4").Value = "'True"
$wsDe.Range("N4").Value = "'"
$wsDe.Range("O4").Value = "'False"
$wsDe.Range("P4").Value = "'"

$wsDe.Hyperlinks.Add(
    $wsDe.Range("A4"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/83bf196a4a2a17b09863015afd929f9aa74d0463/e2e/ca143e9c-70a9-40bc-a684-5517b1eb4cf1.md",
    "",
    "",
    $fileName
) | Out-Null
$wsDe.Hyperlinks.Add(
    $wsDe.Range("I4"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/83bf196a4a2a17b09863015afd929f9aa74d0463/e2e/ca143e9c-70a9-40bc-a684-5517b1eb4cf1.md",
    "",
    "",
    $fileName
) | Out-Null

Write-Host "Handback report row added for ca143e9c-70a9-40bc-a684-5517b1eb4cf1"
